$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add a new day column AB (11-jul) ---------------
$ws1 = $wb.Worksheets.Item(1)

# Header cell AB1 - same header style as the rest of row 1 (copy format
# from AA1 so we don't spawn a brand-new, duplicate cell style).
$ws1.Range("AB1").Value = "11-jul"
$ws1.Range("AA1").Copy()
$ws1.Range("AB1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Hourly prices for 11-jul
$ws1.Range("AB2").Value = 73.12
$ws1.Range("AB3").Value = 62.48
$ws1.Range("AB4").Value = 51.1
$ws1.Range("AB5").Value = 50.61
$ws1.Range("AB6").Value = 44.16
$ws1.Range("AB7").Value = 39.78
$ws1.Range("AB8").Value = 56.01
$ws1.Range("AB9").Value = 65.82
$ws1.Range("AB10").Value = 57.2
$ws1.Range("AB11").Value = 60.05
$ws1.Range("AB12").Value = 50
$ws1.Range("AB13").Value = 28.35
$ws1.Range("AB14").Value = 39.46
$ws1.Range("AB15").Value = 28.11
$ws1.Range("AB16").Value = 32.04
$ws1.Range("AB17").Value = 22.38
$ws1.Range("AB18").Value = 30.7
$ws1.Range("AB19").Value = 47.85
$ws1.Range("AB20").Value = 66.57
$ws1.Range("AB21").Value = 77.7
$ws1.Range("AB22").Value = 71.92
$ws1.Range("AB23").Value = 61.92
$ws1.Range("AB24").Value = 98.93
$ws1.Range("AB25").Value = 89.99

# --- Sheet "Gaz": append 2025-07-09 ------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# A leading apostrophe forces text entry so the ISO-looking date string
# isn't auto-converted into a date serial number; ClearFormats() then
# drops the "quote prefix" cell format Excel would otherwise remember,
# so the cell keeps the workbook's default (unstyled) look.
$ws2.Range("A25").Value = "'2025-07-09"
$ws2.Range("A25").ClearFormats()
$ws2.Range("B25").Value = 33.6

# --- Sheet "CO2": append 2025-07-09 ------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A25").Value = "'2025-07-09"
$ws3.Range("A25").ClearFormats()
$ws3.Range("B25").Value = 69.65
